$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. The document used to carry a stray "_GoBack" bookmark sitting
#    alone in the very first (empty) paragraph. Remove it from there;
#    it gets re-created further down at its real location.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. "CC 12345678" -> "ID 12345678", with the (now orphaned) "_GoBack"
#    bookmark re-inserted right between "ID" and the following space,
#    i.e. splitting the run into "ID" + " 12345678".
# ------------------------------------------------------------------
$d.Content.Find.Execute("CC 12345678", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "ID 12345678", 2)

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "ID 12345678`r") {
        $splitAt = $p.Range.Start + 2
        $bmRange = $d.Range($splitAt, $splitAt)
        $d.Bookmarks.Add("_GoBack", $bmRange)
        break
    }
}
